# Update TPM-normalized NATMI C3-Itgb2 metrics (ligand/receptor expression,
# specificity, and edge-weight columns G:T) for all 16 sending/target cluster
# pairs on the active sheet, per the recomputed-TPM values from the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = [double]"0.555934"
$ws.Cells.Item(2, 8).Value = [double]"1.667802"
$ws.Cells.Item(2, 9).Value = [double]"0.005745252779589096"
$ws.Cells.Item(2, 10).Value = [double]"0.005745252779589094"
$ws.Cells.Item(2, 11).Value = [double]"1"
$ws.Cells.Item(2, 12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(2, 13).Value = [double]"0.06447966666666667"
$ws.Cells.Item(2, 14).Value = [double]"0.193439"
$ws.Cells.Item(2, 15).Value = [double]"0.001101138907643723"
$ws.Cells.Item(2, 16).Value = [double]"0.001101138907643722"
$ws.Cells.Item(2, 17).Value = [double]"0.03584643900866667"
$ws.Cells.Item(2, 18).Value = [double]"0.322617951078"
$ws.Cells.Item(2, 19).Value = [double]"6.326321369853797E-06"
$ws.Cells.Item(2, 20).Value = [double]"6.326321369853795E-06"
$ws.Cells.Item(3, 7).Value = [double]"0.555934"
$ws.Cells.Item(3, 8).Value = [double]"1.667802"
$ws.Cells.Item(3, 9).Value = [double]"0.005745252779589096"
$ws.Cells.Item(3, 10).Value = [double]"0.005745252779589094"
$ws.Cells.Item(3, 15).Value = [double]"0.00657695954769643"
$ws.Cells.Item(3, 16).Value = [double]"0.006576959547696431"
$ws.Cells.Item(3, 17).Value = [double]"0.2141061201746667"
$ws.Cells.Item(3, 18).Value = [double]"1.926955081572"
$ws.Cells.Item(3, 19).Value = [double]"3.778629512264796E-05"
$ws.Cells.Item(3, 20).Value = [double]"3.778629512264795E-05"
$ws.Cells.Item(4, 7).Value = [double]"0.555934"
$ws.Cells.Item(4, 8).Value = [double]"1.667802"
$ws.Cells.Item(4, 9).Value = [double]"0.005745252779589096"
$ws.Cells.Item(4, 10).Value = [double]"0.005745252779589094"
$ws.Cells.Item(4, 11).Value = [double]"1"
$ws.Cells.Item(4, 12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(4, 13).Value = [double]"0.008175"
$ws.Cells.Item(4, 14).Value = [double]"0.024525"
$ws.Cells.Item(4, 15).Value = [double]"0.0001396069650378791"
$ws.Cells.Item(4, 16).Value = [double]"0.0001396069650378791"
$ws.Cells.Item(4, 17).Value = [double]"0.00454476045"
$ws.Cells.Item(4, 18).Value = [double]"0.04090284405"
$ws.Cells.Item(4, 19).Value = [double]"8.020773039338725E-07"
$ws.Cells.Item(4, 20).Value = [double]"8.020773039338722E-07"
$ws.Cells.Item(5, 7).Value = [double]"0.555934"
$ws.Cells.Item(5, 8).Value = [double]"1.667802"
$ws.Cells.Item(5, 9).Value = [double]"0.005745252779589096"
$ws.Cells.Item(5, 10).Value = [double]"0.005745252779589094"
$ws.Cells.Item(5, 13).Value = [double]"58.099467"
$ws.Cells.Item(5, 14).Value = [double]"174.298401"
$ws.Cells.Item(5, 15).Value = [double]"0.992182294579622"
$ws.Cells.Item(5, 16).Value = [double]"0.992182294579622"
$ws.Cells.Item(5, 17).Value = [double]"32.299469087178"
$ws.Cells.Item(5, 18).Value = [double]"290.695221784602"
$ws.Cells.Item(5, 19).Value = [double]"0.00570033808579266"
$ws.Cells.Item(5, 20).Value = [double]"0.005700338085792659"
$ws.Cells.Item(6, 9).Value = [double]"0.823525905561055"
$ws.Cells.Item(6, 10).Value = [double]"0.823525905561055"
$ws.Cells.Item(6, 11).Value = [double]"1"
$ws.Cells.Item(6, 12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(6, 13).Value = [double]"0.06447966666666667"
$ws.Cells.Item(6, 14).Value = [double]"0.193439"
$ws.Cells.Item(6, 15).Value = [double]"0.001101138907643723"
$ws.Cells.Item(6, 16).Value = [double]"0.001101138907643722"
$ws.Cells.Item(6, 17).Value = [double]"5.138237128682556"
$ws.Cells.Item(6, 18).Value = [double]"46.24413415814301"
$ws.Cells.Item(6, 19).Value = [double]"0.0009068164160658075"
$ws.Cells.Item(6, 20).Value = [double]"0.0009068164160658073"
$ws.Cells.Item(7, 9).Value = [double]"0.823525905561055"
$ws.Cells.Item(7, 10).Value = [double]"0.823525905561055"
$ws.Cells.Item(7, 15).Value = [double]"0.00657695954769643"
$ws.Cells.Item(7, 16).Value = [double]"0.006576959547696431"
$ws.Cells.Item(7, 19).Value = [double]"0.005416296567355129"
$ws.Cells.Item(7, 20).Value = [double]"0.00541629656735513"
$ws.Cells.Item(8, 9).Value = [double]"0.823525905561055"
$ws.Cells.Item(8, 10).Value = [double]"0.823525905561055"
$ws.Cells.Item(8, 11).Value = [double]"1"
$ws.Cells.Item(8, 12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(8, 13).Value = [double]"0.008175"
$ws.Cells.Item(8, 14).Value = [double]"0.024525"
$ws.Cells.Item(8, 15).Value = [double]"0.0001396069650378791"
$ws.Cells.Item(8, 16).Value = [double]"0.0001396069650378791"
$ws.Cells.Item(8, 17).Value = [double]"0.651447048325"
$ws.Cells.Item(8, 18).Value = [double]"5.863023434925001"
$ws.Cells.Item(8, 19).Value = [double]"0.0001149699523054499"
$ws.Cells.Item(8, 20).Value = [double]"0.0001149699523054499"
$ws.Cells.Item(9, 9).Value = [double]"0.823525905561055"
$ws.Cells.Item(9, 10).Value = [double]"0.823525905561055"
$ws.Cells.Item(9, 13).Value = [double]"58.099467"
$ws.Cells.Item(9, 14).Value = [double]"174.298401"
$ws.Cells.Item(9, 15).Value = [double]"0.992182294579622"
$ws.Cells.Item(9, 16).Value = [double]"0.992182294579622"
$ws.Cells.Item(9, 17).Value = [double]"4629.813613015994"
$ws.Cells.Item(9, 18).Value = [double]"41668.32251714394"
$ws.Cells.Item(9, 19).Value = [double]"0.8170878226253286"
$ws.Cells.Item(9, 20).Value = [double]"0.8170878226253286"
$ws.Cells.Item(10, 7).Value = [double]"0.3446996666666666"
$ws.Cells.Item(10, 8).Value = [double]"1.034099"
$ws.Cells.Item(10, 9).Value = [double]"0.003562269474506148"
$ws.Cells.Item(10, 10).Value = [double]"0.003562269474506148"
$ws.Cells.Item(10, 11).Value = [double]"1"
$ws.Cells.Item(10, 12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(10, 13).Value = [double]"0.06447966666666667"
$ws.Cells.Item(10, 14).Value = [double]"0.193439"
$ws.Cells.Item(10, 15).Value = [double]"0.001101138907643723"
$ws.Cells.Item(10, 16).Value = [double]"0.001101138907643722"
$ws.Cells.Item(10, 17).Value = [double]"0.02222611960677778"
$ws.Cells.Item(10, 18).Value = [double]"0.200035076461"
$ws.Cells.Item(10, 19).Value = [double]"3.922553517890278E-06"
$ws.Cells.Item(10, 20).Value = [double]"3.922553517890276E-06"
$ws.Cells.Item(11, 7).Value = [double]"0.3446996666666666"
$ws.Cells.Item(11, 8).Value = [double]"1.034099"
$ws.Cells.Item(11, 9).Value = [double]"0.003562269474506148"
$ws.Cells.Item(11, 10).Value = [double]"0.003562269474506148"
$ws.Cells.Item(11, 15).Value = [double]"0.00657695954769643"
$ws.Cells.Item(11, 16).Value = [double]"0.006576959547696431"
$ws.Cells.Item(11, 17).Value = [double]"0.1327537230237778"
$ws.Cells.Item(11, 18).Value = [double]"1.194783507214"
$ws.Cells.Item(11, 19).Value = [double]"2.342890223182076E-05"
$ws.Cells.Item(11, 20).Value = [double]"2.342890223182076E-05"
$ws.Cells.Item(12, 7).Value = [double]"0.3446996666666666"
$ws.Cells.Item(12, 8).Value = [double]"1.034099"
$ws.Cells.Item(12, 9).Value = [double]"0.003562269474506148"
$ws.Cells.Item(12, 10).Value = [double]"0.003562269474506148"
$ws.Cells.Item(12, 11).Value = [double]"1"
$ws.Cells.Item(12, 12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(12, 13).Value = [double]"0.008175"
$ws.Cells.Item(12, 14).Value = [double]"0.024525"
$ws.Cells.Item(12, 15).Value = [double]"0.0001396069650378791"
$ws.Cells.Item(12, 16).Value = [double]"0.0001396069650378791"
$ws.Cells.Item(12, 17).Value = [double]"0.002817919775"
$ws.Cells.Item(12, 18).Value = [double]"0.025361277975"
$ws.Cells.Item(12, 19).Value = [double]"4.973176299828837E-07"
$ws.Cells.Item(12, 20).Value = [double]"4.973176299828837E-07"
$ws.Cells.Item(13, 7).Value = [double]"0.3446996666666666"
$ws.Cells.Item(13, 8).Value = [double]"1.034099"
$ws.Cells.Item(13, 9).Value = [double]"0.003562269474506148"
$ws.Cells.Item(13, 10).Value = [double]"0.003562269474506148"
$ws.Cells.Item(13, 13).Value = [double]"58.099467"
$ws.Cells.Item(13, 14).Value = [double]"174.298401"
$ws.Cells.Item(13, 15).Value = [double]"0.992182294579622"
$ws.Cells.Item(13, 16).Value = [double]"0.992182294579622"
$ws.Cells.Item(13, 17).Value = [double]"20.026866908411"
$ws.Cells.Item(13, 18).Value = [double]"180.241802175699"
$ws.Cells.Item(13, 19).Value = [double]"0.003534420701126455"
$ws.Cells.Item(13, 20).Value = [double]"0.003534420701126454"
$ws.Cells.Item(14, 7).Value = [double]"16.17571666666667"
$ws.Cells.Item(14, 8).Value = [double]"48.52715"
$ws.Cells.Item(14, 9).Value = [double]"0.1671665721848498"
$ws.Cells.Item(14, 10).Value = [double]"0.1671665721848498"
$ws.Cells.Item(14, 11).Value = [double]"1"
$ws.Cells.Item(14, 12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(14, 13).Value = [double]"0.06447966666666667"
$ws.Cells.Item(14, 14).Value = [double]"0.193439"
$ws.Cells.Item(14, 15).Value = [double]"0.001101138907643723"
$ws.Cells.Item(14, 16).Value = [double]"0.001101138907643722"
$ws.Cells.Item(14, 17).Value = [double]"1.043004818761111"
$ws.Cells.Item(14, 18).Value = [double]"9.38704336885"
$ws.Cells.Item(14, 19).Value = [double]"0.0001840736166901711"
$ws.Cells.Item(14, 20).Value = [double]"0.000184073616690171"
$ws.Cells.Item(15, 7).Value = [double]"16.17571666666667"
$ws.Cells.Item(15, 8).Value = [double]"48.52715"
$ws.Cells.Item(15, 9).Value = [double]"0.1671665721848498"
$ws.Cells.Item(15, 10).Value = [double]"0.1671665721848498"
$ws.Cells.Item(15, 15).Value = [double]"0.00657695954769643"
$ws.Cells.Item(15, 16).Value = [double]"0.006576959547696431"
$ws.Cells.Item(15, 17).Value = [double]"6.229732192211111"
$ws.Cells.Item(15, 18).Value = [double]"56.0675897299"
$ws.Cells.Item(15, 19).Value = [double]"0.001099447782986833"
$ws.Cells.Item(15, 20).Value = [double]"0.001099447782986833"
$ws.Cells.Item(16, 7).Value = [double]"16.17571666666667"
$ws.Cells.Item(16, 8).Value = [double]"48.52715"
$ws.Cells.Item(16, 9).Value = [double]"0.1671665721848498"
$ws.Cells.Item(16, 10).Value = [double]"0.1671665721848498"
$ws.Cells.Item(16, 11).Value = [double]"1"
$ws.Cells.Item(16, 12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(16, 13).Value = [double]"0.008175"
$ws.Cells.Item(16, 14).Value = [double]"0.024525"
$ws.Cells.Item(16, 15).Value = [double]"0.0001396069650378791"
$ws.Cells.Item(16, 16).Value = [double]"0.0001396069650378791"
$ws.Cells.Item(16, 17).Value = [double]"0.13223648375"
$ws.Cells.Item(16, 18).Value = [double]"1.19012835375"
$ws.Cells.Item(16, 19).Value = [double]"2.333761779851242E-05"
$ws.Cells.Item(16, 20).Value = [double]"2.333761779851242E-05"
$ws.Cells.Item(17, 7).Value = [double]"16.17571666666667"
$ws.Cells.Item(17, 8).Value = [double]"48.52715"
$ws.Cells.Item(17, 9).Value = [double]"0.1671665721848498"
$ws.Cells.Item(17, 10).Value = [double]"0.1671665721848498"
$ws.Cells.Item(17, 13).Value = [double]"58.099467"
$ws.Cells.Item(17, 14).Value = [double]"174.298401"
$ws.Cells.Item(17, 15).Value = [double]"0.992182294579622"
$ws.Cells.Item(17, 16).Value = [double]"0.992182294579622"
$ws.Cells.Item(17, 17).Value = [double]"939.8005166763501"
$ws.Cells.Item(17, 18).Value = [double]"8458.204650087151"
$ws.Cells.Item(17, 19).Value = [double]"0.1658597131673743"
$ws.Cells.Item(17, 20).Value = [double]"0.1658597131673743"

Write-Output "Updated 190 cells (columns G:T, rows 2-17) with recomputed TPM values."
